$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (G) values - regenerated to use K instead of Strike#
$values = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 1
    6 = 1
    7 = 4
    8 = 2
    9 = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 2
    16 = 2
    17 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
